# Estadisticos Segundo Parcial 26 Mayo
#
# "Rescatables" sheet: four students (ALVAREZ, COLMENARES, CRESCENCIO, PEREZ)
# are sorted by number of subjects still failed (column G, "Reprobadas"),
# descending. After the second-partial update ALVAREZ's failure count drops
# from 2 to 1, so that record now sorts to the bottom of the block and the
# other three rows shift up by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Row 19 <- old row 20 (COLMENARES)
$ws.Range("A19").Value() = 22330051920424
$ws.Range("B19").Value() = "COLMENARES"
$ws.Range("C19").Value() = "MARTINEZ"
$ws.Range("D19").Value() = "JULIO EDUARDO"
$ws.Range("E19").Value() = "Ciencias sociales III"
$ws.Range("F19").Value() = "4BLCM"
$ws.Range("G19").Value() = 2

# Row 20 <- old row 21 (CRESCENCIO)
$ws.Range("A20").Value() = 22330051920033
$ws.Range("B20").Value() = "CRESCENCIO"
$ws.Range("C20").Value() = "DIAZ"
$ws.Range("D20").Value() = "DIEGO ARMANDO"
$ws.Range("E20").Value() = "TEMAS DE FILOSOFÍA"
$ws.Range("F20").Value() = "6BEM"
$ws.Range("G20").Value() = 2

# Row 21 <- old row 22 (PEREZ)
$ws.Range("A21").Value() = 22330051920045
$ws.Range("B21").Value() = "PEREZ"
$ws.Range("C21").Value() = "ROMERO"
$ws.Range("D21").Value() = "JULIAN DAVID"
$ws.Range("E21").Value() = "TEMAS DE FILOSOFÍA"
$ws.Range("F21").Value() = "6BEM"
$ws.Range("G21").Value() = 2

# Row 22 <- old row 19 (ALVAREZ), Reprobadas now 1 instead of 2
$ws.Range("A22").Value() = 22330061460232
$ws.Range("B22").Value() = "ALVAREZ"
$ws.Range("C22").Value() = "VOTE"
$ws.Range("D22").Value() = "CAMILO"
$ws.Range("E22").Value() = "Ciencias sociales III"
$ws.Range("F22").Value() = "4BLCM"
$ws.Range("G22").Value() = 1
